$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "92.132.98"
$ws.Cells.Item(2, 5).Value = "  +1.26%  "

$ws.Cells.Item(3, 4).Value = "3.114.13"
$ws.Cells.Item(3, 5).Value = "  -1.22%  "

$ws.Cells.Item(4, 5).Value = "  +0.09%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "236.78"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.18%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "615.11"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.74%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "1.11"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -2.01%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.387"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +4.16%  "

$ws.Cells.Item(9, 5).Value = "  -0.03%  "

$ws.Cells.Item(10, 4).Value = "3.110.02"
$ws.Cells.Item(10, 5).Value = "  -1.45%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.775"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +4.97%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.199"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -1.92%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000246"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.13%  "

$ws.Cells.Item(14, 4).Value = "92.131.70"
$ws.Cells.Item(14, 5).Value = "  +1.72%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "34.05"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -3.23%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "5.43"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -4.05%  "

$ws.Cells.Item(17, 4).Value = "3.705.24"
$ws.Cells.Item(17, 5).Value = "  -0.77%  "

$ws.Cells.Item(18, 4).Value = "3.096.33"
$ws.Cells.Item(18, 5).Value = "  -0.47%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "3.77"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +1.54%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "14.55"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -3.37%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "5.85"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -1.91%  "

$ws.Cells.Item(22, 2).Value = "BitcoinCash"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "443.72"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -2.31%  "

$ws.Cells.Item(23, 2).Value = "Uniswap"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "9.26"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.76%  "

$ws.Cells.Item(24, 5).Value = "  -3.67%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "5.70"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -5.93%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "86.38"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.80%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "11.64"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -4.24%  "

$ws.Cells.Item(28, 5).Value = "  -0.05%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "0.132"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -6.70%  "

$ws.Cells.Item(30, 2).Value = "Cronos"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "0.182"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +8.28%  "

$ws.Cells.Item(31, 2).Value = "Stellar"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "0.235"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -1.23%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "9.10"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -6.18%  "

$ws.Cells.Item(33, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.969"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -2.36%  "

$ws.Cells.Item(34, 2).Value = "RenderToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "7.83"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.59%  "

$ws.Cells.Item(35, 2).Value = "Kaspa"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.158"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -7.97%  "

$ws.Cells.Item(36, 2).Value = "EthereumClassic"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "26.00"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -2.28%  "

$ws.Cells.Item(37, 2).Value = "PancakeSwap"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "1.89"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -4.28%  "

$ws.Cells.Item(38, 2).Value = "MantraDAO"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "3.88"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.67%  "

$ws.Cells.Item(39, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "23.83"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +7.81%  "

$ws.Cells.Item(40, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.439"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.68%  "

$ws.Cells.Item(41, 2).Value = "Fetch.AI"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "1.28"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -5.14%  "

$ws.Cells.Item(42, 2).Value = "Bittensor"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "470.45"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -6.42%  "

$ws.Cells.Item(43, 2).Value = "dogwifhat"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "3.29"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -5.18%  "

$ws.Cells.Item(44, 2).Value = "USDe"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.01%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "161.31"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +3.70%  "

$ws.Cells.Item(46, 2).Value = "ARBITRUM"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.686"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -5.06%  "

$ws.Cells.Item(47, 2).Value = "Stacks"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "1.85"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -3.68%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "1.35"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -2.68%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "0.0330"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.40%  "

$ws.Cells.Item(50, 2).Value = "OKB"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "44.02"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.31%  "

$ws.Cells.Item(51, 2).Value = "Filecoin"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "4.40"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -3.40%  "
